$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.476.63'
$ws.Range('E2').Value = '  +3.40%  '
$ws.Range('D3').Value = '1.601.84'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.64'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.65'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0911'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.97%  '
$ws.Range('D13').Value = '1.833.39'
$ws.Range('E13').Value = '  +3.16%  '
$ws.Range('D14').Value = '1.608.51'
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range('D15').Value = '29.506.04'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.534'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.73'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.56%  '
$ws.Range('D21').Value = '0.0₃0692'
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.48%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.48%  '
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('D34').Value = '1.425.51'
$ws.Range('E34').Value = '  +2.14%  '
$ws.Range('E35').Value = '  +2.60%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.81'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.534'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '53.58'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +21.81%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.792'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0471'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.21%  '
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('D49').Value = '1.743.99'
$ws.Range('E49').Value = '  +3.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.837'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.64%  '
